$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark that used to sit at the end of the
#    paragraph ending in "...SLA-preferences incompatibilities."
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Rewrite the "data integration process" sentence (second challenge
#    paragraph) with the new, expanded wording.
# ------------------------------------------------------------------
$old = "he data integration process includes looking for services that can be used as data providers, and for services required in order (i) to retrieve the data; (ii) to build an integrated result; and (iii) to deliver it to the user considering the user quality requirements, her context and resources consumption."
$new = "he data integration process includes (i) looking up services that can be used as data providers, and for services required to process retrieved data and build an integrated result; (ii) performing data retrieval, processing and integration and (iii) deliver results to the user considering her preferences (quality requirements, context and resources consumption). The integrated SLA can guide services filtering in the look up phase; it can help to control the amounts of data to retrieve and process according to consumption rights depending on the user subscription to the participating cloud providers and how to deliver data considering the user" + [char]0x2019 + "s context."

$d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# ------------------------------------------------------------------
# 3. The whole-text replacement above collapses formatting to a single
#    run; restore italics on the "(iii) " marker, which keeps its
#    original italic formatting in the target document.
# ------------------------------------------------------------------
$p8 = $d.Paragraphs.Item(8)
$scope = $d.Range($p8.Range.Start, $p8.Range.End)
$found = $scope.Find.Execute("(iii) ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $scope.Font.Italic = $true
}

# ------------------------------------------------------------------
# 4. Re-create the "_GoBack" bookmark at the very end of this paragraph
#    (immediately after "...resources consumption." and before the
#    paragraph mark), matching its new location in the target document.
#
#    A truly collapsed range placed exactly at a paragraph-end position
#    trips up Bookmarks.Add in this environment, so we work around it
#    by bookmarking a temporary placeholder character and then deleting
#    that character (the bookmark collapses back down and stays put).
# ------------------------------------------------------------------
$p8 = $d.Paragraphs.Item(8)
$endPos = $p8.Range.End - 1

$beforeChar = $d.Range($endPos - 1, $endPos)
$beforeChar.InsertAfter("X")

$placeholder = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $placeholder) | Out-Null

$placeholder = $d.Range($endPos, $endPos + 1)
$placeholder.Delete()
